$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold + border) from an existing header cell (e.g. E1) to F1:H1
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()
$excel.CutCopyMode = $false

# Fill rows 2-21 in columns F, G, H with boolean FALSE
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
